# Applies the authored changes to the IOA data workbook:
#  1. Column W (timeframe), row 2: "Short Term" -> "Medium Term (5-10 years to generate return)"
#  2. Column AQ (impact_dimensions), rows 2-10: drop the trailing
#     "Contribution: \n\nHow much: " placeholder lines, leaving a single
#     trailing blank line.
#  3. Remove the trailing "image" column (BF) that listed embedded picture
#     filenames - the column and its header are deleted outright.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the timeframe value for the first data row only.
$ws.Cells.Item(2, 23).Value = "Medium Term (5-10 years to generate return)"

# 2) Strip the "Contribution:/How much:" placeholder suffix from the
#    impact_dimensions column (AQ = column 43) for every data row (2-10).
$suffix = "`n`nContribution: `n`nHow much: "
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 43)
    $value = $cell.Value2
    if ($value -ne $null -and $value.EndsWith($suffix)) {
        $trimmed = $value.Substring(0, $value.Length - $suffix.Length)
        $cell.Value = $trimmed + "`n "
    }
}

# 3) Drop the "image" column (BF) entirely - header + the per-row filenames.
$ws.Columns("BF").Delete()
